# Regenerate merged AHB files
# - Rename the "_old" / "_new" header-row columns to "_FV2410" / "_FV2504"
# - Turn the data range into an Excel Table ("Table1") with AutoFilter
# - Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (A1:J1 = "*_old" -> "*_FV2410", L1:U1 = "*_new" -> "*_FV2504") ---
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2410"
}
# Column K1 stays "diff"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, 12 + $i).Value = $baseNames[$i] + "_FV2504"
}

# --- 2. Convert the used range A1:U56 into a Table with an AutoFilter ---
$tableRange = $ws.Range("A1:U56")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# --- 3. Freeze panes at row 2 (keep header row 1 visible while scrolling) ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "AHB table regenerated: headers renamed, Table1 created, panes frozen."
